# Auto-generated Excel COM-interop edit script
# Applies the numeric updates described by the commit diff to the
# "Anima_Profits" profit-tracking tables across each crafting-job sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 333.33334
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -430
$ws.Range("N12").Value = -540
# Row 113
$ws.Range("H113").Value = 3197.6667
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 3515.8
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 3515.8
$ws.Range("M113").Value = 454
$ws.Range("N113").Value = -10023.8
# Row 125
$ws.Range("H125").Value = 1764.9231
$ws.Range("I125").Value = 750
$ws.Range("J125").Value = 1949.4546
$ws.Range("K125").Value = 6750
$ws.Range("L125").Value = 17545.0914
$ws.Range("M125").Value = -4290
$ws.Range("N125").Value = -22465.0914
# Row 129
$ws.Range("H129").Value = 1006.68854
$ws.Range("I129").Value = 642.6667
$ws.Range("J129").Value = 1069.6923
$ws.Range("K129").Value = 1928.0001
$ws.Range("L129").Value = 3209.0769
$ws.Range("M129").Value = 3071.9999
$ws.Range("N129").Value = -13209.0769
# Row 137
$ws.Range("H137").Value = 2897.1
$ws.Range("I137").Value = 2372.1516
$ws.Range("K137").Value = 7116.4548
$ws.Range("M137").Value = -4566.4548
# Row 141
$ws.Range("H141").Value = 3164.6511
$ws.Range("I141").Value = 947.9091
$ws.Range("J141").Value = 10479.9
$ws.Range("K141").Value = 2843.7273
$ws.Range("L141").Value = 31439.7
$ws.Range("M141").Value = 2336.2727
$ws.Range("N141").Value = -41799.7

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 555490
$ws.Range("I32").Value = 619607.25
$ws.Range("J32").Value = 24232.572
$ws.Range("K32").Value = 619607.25
$ws.Range("L32").Value = 24232.572
$ws.Range("M32").Value = -619320.25
$ws.Range("N32").Value = -24806.572
# Row 61
$ws.Range("H61").Value = 2035.463
$ws.Range("I61").Value = 1630.4359
$ws.Range("J61").Value = 3088.5334
$ws.Range("K61").Value = 1630.4359
$ws.Range("L61").Value = 3088.5334
$ws.Range("M61").Value = -1418.4359
$ws.Range("N61").Value = -3512.5334
# Row 63
$ws.Range("H63").Value = 6722.5557
$ws.Range("I63").Value = 5333.3335
$ws.Range("J63").Value = 7417.1665
$ws.Range("K63").Value = 5333.3335
$ws.Range("L63").Value = 7417.1665
$ws.Range("M63").Value = -4647.3335
$ws.Range("N63").Value = -8789.166499999999
# Row 66
$ws.Range("H66").Value = 6722.5557
$ws.Range("I66").Value = 5333.3335
$ws.Range("J66").Value = 7417.1665
$ws.Range("K66").Value = 26666.6675
$ws.Range("L66").Value = 37085.8325
$ws.Range("M66").Value = -23234.6675
$ws.Range("N66").Value = -43949.8325
# Row 74
$ws.Range("H74").Value = 1078.2954
$ws.Range("I74").Value = 728.24243
$ws.Range("J74").Value = 2128.4546
$ws.Range("K74").Value = 728.24243
$ws.Range("L74").Value = 2128.4546
$ws.Range("M74").Value = 145.75757
$ws.Range("N74").Value = -3876.4546
# Row 77
$ws.Range("H77").Value = 1078.2954
$ws.Range("I77").Value = 728.24243
$ws.Range("J77").Value = 2128.4546
$ws.Range("K77").Value = 3641.21215
$ws.Range("L77").Value = 10642.273
$ws.Range("M77").Value = 726.7878499999997
$ws.Range("N77").Value = -19378.273
# Row 97
$ws.Range("H97").Value = 972.8333
$ws.Range("I97").Value = 878.5714
$ws.Range("J97").Value = 1302.75
$ws.Range("K97").Value = 878.5714
$ws.Range("L97").Value = 1302.75
$ws.Range("M97").Value = -382.5714
$ws.Range("N97").Value = -2294.75
# Row 122
$ws.Range("H122").Value = 112477.78
$ws.Range("I122").Value = 167616.67
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 502850.01
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -500400.01
$ws.Range("N122").Value = -11500
# Row 128
$ws.Range("H128").Value = 40543
$ws.Range("J128").Value = 40543
$ws.Range("L128").Value = 40543
$ws.Range("N128").Value = -50503
# Row 136
$ws.Range("H136").Value = 2035.463
$ws.Range("I136").Value = 1630.4359
$ws.Range("J136").Value = 3088.5334
$ws.Range("K136").Value = 4891.307699999999
$ws.Range("L136").Value = 9265.600199999999
$ws.Range("M136").Value = -2341.307699999999
$ws.Range("N136").Value = -14365.6002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1823.3721
$ws.Range("I134").Value = 1475.5294
$ws.Range("J134").Value = 3137.4443
$ws.Range("K134").Value = 4426.5882
$ws.Range("L134").Value = 9412.332900000001
$ws.Range("M134").Value = -1891.5882
$ws.Range("N134").Value = -14482.3329

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 576.4706
$ws.Range("J22").Value = 800
$ws.Range("L22").Value = 800
$ws.Range("N22").Value = -1500
# Row 31
$ws.Range("H31").Value = 5171.712
$ws.Range("I31").Value = 1379.2667
$ws.Range("J31").Value = 9094.931
$ws.Range("K31").Value = 1379.2667
$ws.Range("L31").Value = 9094.931
$ws.Range("M31").Value = -1084.2667
$ws.Range("N31").Value = -9684.931
# Row 34
$ws.Range("H34").Value = 5171.712
$ws.Range("I34").Value = 1379.2667
$ws.Range("J34").Value = 9094.931
$ws.Range("K34").Value = 1379.2667
$ws.Range("L34").Value = 9094.931
$ws.Range("M34").Value = -1177.2667
$ws.Range("N34").Value = -9498.931
# Row 99
$ws.Range("H99").Value = 1975.68
$ws.Range("I99").Value = 1749.5
$ws.Range("K99").Value = 1749.5
$ws.Range("M99").Value = -251.5
# Row 126
$ws.Range("H126").Value = 1975.68
$ws.Range("I126").Value = 1749.5
$ws.Range("K126").Value = 5248.5
$ws.Range("M126").Value = -2778.5
# Row 132
$ws.Range("H132").Value = 5210433.5
$ws.Range("I132").Value = 1552.6842
$ws.Range("K132").Value = 4658.0526
$ws.Range("M132").Value = -2128.0526

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 832.2353000000001
$ws.Range("I5").Value = 578.1818
$ws.Range("J5").Value = 1298
$ws.Range("K5").Value = 1734.5454
$ws.Range("L5").Value = 3894
$ws.Range("M5").Value = -1622.5454
$ws.Range("N5").Value = -4118
# Row 31
$ws.Range("H31").Value = 1661.625
$ws.Range("J31").Value = 1572.4
$ws.Range("L31").Value = 4717.200000000001
$ws.Range("N31").Value = -5293.200000000001
# Row 132
$ws.Range("H132").Value = 2362.152
$ws.Range("I132").Value = 2323.0908
$ws.Range("J132").Value = 2374.4285
$ws.Range("K132").Value = 20907.8172
$ws.Range("L132").Value = 21369.8565
$ws.Range("M132").Value = -18377.8172
$ws.Range("N132").Value = -26429.8565
# Row 134
$ws.Range("H134").Value = 6741.6665
$ws.Range("I134").Value = 2969.8572
$ws.Range("J134").Value = 7652.1035
$ws.Range("K134").Value = 8909.571599999999
$ws.Range("L134").Value = 22956.3105
$ws.Range("M134").Value = -3839.571599999999
$ws.Range("N134").Value = -33096.3105
# Row 135
$ws.Range("H135").Value = 832.2353000000001
$ws.Range("I135").Value = 578.1818
$ws.Range("J135").Value = 1298
$ws.Range("K135").Value = 5203.6362
$ws.Range("L135").Value = 11682
$ws.Range("M135").Value = -2668.6362
$ws.Range("N135").Value = -16752

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2751.111
$ws.Range("I97").Value = 2870
$ws.Range("J97").Value = 1800
$ws.Range("K97").Value = 2870
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -2374
$ws.Range("N97").Value = -2792
# Row 102
$ws.Range("H102").Value = 1975.7273
$ws.Range("I102").Value = 1779.125
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1779.125
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -157.125
$ws.Range("N102").Value = -5744
# Row 132
$ws.Range("H132").Value = 2141.446
$ws.Range("I132").Value = 1803.3455
$ws.Range("J132").Value = 4001
$ws.Range("K132").Value = 5410.0365
$ws.Range("L132").Value = 12003
$ws.Range("M132").Value = -2880.0365
$ws.Range("N132").Value = -17063
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4000.611
$ws.Range("I7").Value = 2615.5715
$ws.Range("J7").Value = 4882
$ws.Range("K7").Value = 2615.5715
$ws.Range("L7").Value = 4882
$ws.Range("M7").Value = -2503.5715
$ws.Range("N7").Value = -5106
# Row 40
$ws.Range("H40").Value = 113704.22
$ws.Range("I40").Value = 145334
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 145334
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -145198
$ws.Range("N40").Value = -3272
# Row 126
$ws.Range("H126").Value = 4000.611
$ws.Range("I126").Value = 2615.5715
$ws.Range("J126").Value = 4882
$ws.Range("K126").Value = 7846.7145
$ws.Range("L126").Value = 14646
$ws.Range("M126").Value = -5376.7145
$ws.Range("N126").Value = -19586
# Row 132
$ws.Range("H132").Value = 1690.5883
$ws.Range("I132").Value = 1192.075
$ws.Range("J132").Value = 3503.3635
$ws.Range("K132").Value = 3576.225
$ws.Range("L132").Value = 10510.0905
$ws.Range("M132").Value = -1046.225
$ws.Range("N132").Value = -15570.0905

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2202.5186
$ws.Range("I122").Value = 2127.9
$ws.Range("K122").Value = 6383.700000000001
$ws.Range("M122").Value = -3933.700000000001
# Row 126
$ws.Range("H126").Value = 5650.625
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 6534.1665
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 19602.4995
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -24542.4995
# Row 132
$ws.Range("H132").Value = 5210830
$ws.Range("I132").Value = 2807.111
$ws.Range("J132").Value = 11906859
$ws.Range("K132").Value = 8421.332999999999
$ws.Range("L132").Value = 35720577
$ws.Range("M132").Value = -5891.332999999999
$ws.Range("N132").Value = -35725637

Write-Host "Applied Anima_Profits scheduled runner updates across all sheets."